$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.384.22'

$ws.Range('E2').Value = '  +1.77%  '

$ws.Range('D3').Value = '1.850.93'

$ws.Range('E3').Value = '  +0.96%  '

$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.0000'
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '245.03'
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.6919'
$c.Style = "Normal"

$ws.Range('E6').Value = '  +0.25%  '

$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"

$ws.Range('E7').Value = '  +0.03%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3054'
$c.Style = "Normal"

$ws.Range('E8').Value = '  +0.18%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.07630'
$c.Style = "Normal"

$ws.Range('E9').Value = '  -0.95%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '23.44'
$c.Style = "Normal"

$ws.Range('E10').Value = '  +0.10%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07726'
$c.Style = "Normal"

$ws.Range('E11').Value = '  -1.07%  '

$ws.Range('D12').Value = '1.853.25'

$ws.Range('E12').Value = '  +1.02%  '

$ws.Range('E13').Value = '  +0.75%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.6916'
$c.Style = "Normal"

$ws.Range('E14').Value = '  +1.32%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '90.49'
$c.Style = "Normal"

$ws.Range('E15').Value = '  -1.69%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '6.289'
$c.Style = "Normal"

$ws.Range('E16').Value = '  -2.33%  '

$ws.Range('D17').Value = '29.396.24'

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '0.000008250'
$c.Style = "Normal"

$ws.Range('E18').Value = '  -0.70%  '

$ws.Range('D19').Value = '2.095.96'

$ws.Range('E19').Value = '  +1.10%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '235.91'
$c.Style = "Normal"

$ws.Range('E20').Value = '  -2.78%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '12.68'
$c.Style = "Normal"

$ws.Range('E21').Value = '  -0.32%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '7.665'
$c.Style = "Normal"

$ws.Range('E23').Value = '  +2.92%  '

$ws.Range('E24').Value = '  +0.09%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.1470'
$c.Style = "Normal"

$ws.Range('E25').Value = '  -0.40%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '8.935'
$c.Style = "Normal"

$ws.Range('E26').Value = '  +1.57%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '160.00'
$c.Style = "Normal"

$ws.Range('E27').Value = '  +1.16%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '18.18'
$c.Style = "Normal"

$ws.Range('E28').Value = '  -0.36%  '

$ws.Range('E29').Value = '  -1.00%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '4.252'
$c.Style = "Normal"

$ws.Range('E30').Value = '  +0.70%  '

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '4.133'
$c.Style = "Normal"

$ws.Range('E31').Value = '  -0.74%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '1.200'
$c.Style = "Normal"

$ws.Range('E32').Value = '  +0.35%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.05218'
$c.Style = "Normal"

$ws.Range('E33').Value = '  +2.44%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.7752'
$c.Style = "Normal"

$ws.Range('E34').Value = '  -0.33%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.872'
$c.Style = "Normal"

$ws.Range('E35').Value = '  +1.36%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.146'
$c.Style = "Normal"

$ws.Range('E36').Value = '  +0.36%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.678'
$c.Style = "Normal"

$ws.Range('E37').Value = '  -0.52%  '

$ws.Range('D38').Value = '1.312.78'

$ws.Range('E38').Value = '  +6.24%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.01861'
$c.Style = "Normal"

$ws.Range('E39').Value = '  +0.45%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '2.709'
$c.Style = "Normal"

$ws.Range('E40').Value = '  +0.57%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.9403'
$c.Style = "Normal"

$ws.Range('E41').Value = '  -1.34%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '105.73'
$c.Style = "Normal"

$ws.Range('E42').Value = '  -2.75%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '5.774'
$c.Style = "Normal"

$ws.Range('E43').Value = '  -2.06%  '

$ws.Range('E44').Value = '  +0.00%  '

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '9.699'
$c.Style = "Normal"

$ws.Range('E45').Value = '  +1.03%  '

$ws.Range('D46').Value = '1.997.32'

$ws.Range('E46').Value = '  +1.17%  '

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.5216'
$c.Style = "Normal"

$ws.Range('E47').Value = '  +1.16%  '

$ws.Range('B48').Value = 'RenderToken'

$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.775'
$c.Style = "Normal"

$ws.Range('E48').Value = '  +1.55%  '

$ws.Range('B49').Value = 'BabyDogeCoin'

$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.00000000120'
$c.Style = "Normal"

$ws.Range('E49').Value = '  -1.62%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '63.06'
$c.Style = "Normal"

$ws.Range('E50').Value = '  -1.63%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.05950'
$c.Style = "Normal"

$ws.Range('E51').Value = '  +0.97%  '

